$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the A:E values of the paired rows (each "B" label row swaps with the
# following "C" label row), before the F/G columns are removed.
# Row 3 (2005年B) <-> Row 4 (2005年C)
$ws.Range("A3").Value = "2005年C"
$ws.Range("B3").Value = 97.8
$ws.Range("C3").Value = 11.5
$ws.Range("E3").Value = 505.1

$ws.Range("A4").Value = "2005年B"
$ws.Range("B4").Value = 98.7
$ws.Range("C4").Value = 22.2
$ws.Range("E4").Value = 326.1

# Row 7 (2016年B) <-> Row 8 (2016年C)
$ws.Range("A7").Value = "2016年C"
$ws.Range("B7").Value = 99.5
$ws.Range("C7").Value = -16.4
$ws.Range("D7").Value = 1.7
$ws.Range("E7").Value = 2291.84138

$ws.Range("A8").Value = "2016年B"
$ws.Range("B8").Value = 99.7
$ws.Range("C8").Value = -9.300000000000001
$ws.Range("D8").Value = 1.9
$ws.Range("E8").Value = 1502.92436

# Row 11 (2017年B) <-> Row 12 (2017年C)
$ws.Range("A11").Value = "2017年C"
$ws.Range("B11").Value = 96.7
$ws.Range("C11").Value = 30.1
$ws.Range("D11").Value = -3.5
$ws.Range("E11").Value = 2327.44371

$ws.Range("A12").Value = "2017年B"
$ws.Range("B12").Value = 99
$ws.Range("C12").Value = 20.1
$ws.Range("D12").Value = -0.6
$ws.Range("E12").Value = 1627.56517

# Row 15 (2018年B) <-> Row 16 (2018年C)
$ws.Range("A15").Value = "2018年C"
$ws.Range("B15").Value = 99.7
$ws.Range("C15").Value = -1.3
$ws.Range("D15").Value = 0.3
$ws.Range("E15").Value = 2416.1711

$ws.Range("A16").Value = "2018年B"
$ws.Range("B16").Value = 100.1
$ws.Range("C16").Value = -1.4
$ws.Range("D16").Value = 0.7
$ws.Range("E16").Value = 1601.62563

# Row 19 (2019年B) <-> Row 20 (2019年C)
$ws.Range("A19").Value = "2019年C"
$ws.Range("B19").Value = 100.1
$ws.Range("C19").Value = -8.4
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 2488.10479

$ws.Range("A20").Value = "2019年B"
$ws.Range("B20").Value = 99.7
$ws.Range("C20").Value = 6.8
$ws.Range("D20").Value = 0.2
$ws.Range("E20").Value = 1637.68901

# Remove columns F and G entirely (产销率 / 销售量 duplicate columns).
$ws.Range("F1:G21").EntireColumn.Delete()
